$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Underaged"
$ws.Range("C13").Value = 19
$ws.Range("D13").Value = $true
